# Safety measures for folder deletion
# Remove the "Global Tolerance" row (row 3) from the worksheet: select the
# entire row first (as a user would when right-clicking the row header and
# choosing "Delete"), then delete it so all rows below shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Rows(3)
$row.Select()
$row.Delete()
